$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(2, "ECs", "Il1b", "Il1rap", "ECs", 3, 1, 0.197811, 0.593433, 0.001274141111268169, 0.001274141111268169, 3, 1, 2.652098, 7.956294000000001, 0.1783170651171435, 0.1783170651171436, 0.524614157478, 4.721527417302, 0.0002272011035064357, 0.0002272011035064358),
  @(3, "ECs", "Il1b", "Il1rap", "FAPs", 3, 1, 0.197811, 0.593433, 0.001274141111268169, 0.001274141111268169, 3, 1, 5.511188000000001, 16.533564, 0.370551491486672, 0.3705514914866722, 1.090173609468, 9.811562485212, 0.0004721348891449057, 0.000472134889144906),
  @(4, "ECs", "Il1b", "Il1rap", "M2", 3, 1, 0.197811, 0.593433, 0.001274141111268169, 0.001274141111268169, 3, 1, 2.28691, 6.860729999999999, 0.1537632015811809, 0.1537632015811809, 0.4523759540099999, 4.07138358609, 0.0001959160165347972, 0.0001959160165347973),
  @(5, "ECs", "Il1b", "Il1rap", "sCs", 3, 1, 0.197811, 0.593433, 0.001274141111268169, 0.001274141111268169, 3, 1, 4.422738333333334, 13.268215, 0.2973682418150034, 0.2973682418150034, 0.8748662924550001, 7.873796632095001, 0.0003788891020820299, 0.0003788891020820301),
  @(6, "M2", "Il1b", "Il1rap", "ECs", 3, 1, 154.0066273333333, 462.0198820000001, 0.9919881871744044, 0.9919881871744045, 3, 1, 2.652098, 7.956294000000001, 0.1783170651171435, 0.1783170651171436, 408.4406683374787, 3675.966015037309, 0.1768884221678154, 0.1768884221678155),
  @(7, "M2", "Il1b", "Il1rap", "FAPs", 3, 1, 154.0066273333333, 462.0198820000001, 0.9919881871744044, 0.9919881871744045, 3, 1, 5.511188000000001, 16.533564, 0.370551491486672, 0.3705514914866722, 848.7594764799388, 7638.835288319449, 0.3675827022946355, 0.3675827022946357),
  @(8, "M2", "Il1b", "Il1rap", "M2", 3, 1, 154.0066273333333, 462.0198820000001, 0.9919881871744044, 0.9919881871744045, 3, 1, 2.28691, 6.860729999999999, 0.1537632015811809, 0.1537632015811809, 352.1992961148733, 3169.79366503386, 0.1525312795906481, 0.1525312795906482),
  @(9, "M2", "Il1b", "Il1rap", "sCs", 3, 1, 154.0066273333333, 462.0198820000001, 0.9919881871744044, 0.9919881871744045, 3, 1, 4.422738333333334, 13.268215, 0.2973682418150034, 0.2973682418150034, 681.1310142945146, 6130.179128650631, 0.2949857831213051, 0.2949857831213052),
  @(10, "sCs", "Il1b", "Il1rap", "ECs", 1, 0.3333333333333333, 1.046026666666667, 3.13808, 0.006737671714327339, 0.006737671714327339, 3, 1, 2.652098, 7.956294000000001, 0.1783170651171435, 0.1783170651171436, 2.774165230613333, 24.96748707552, 0.001201441845821644, 0.001201441845821644),
  @(11, "sCs", "Il1b", "Il1rap", "FAPs", 1, 0.3333333333333333, 1.046026666666667, 3.13808, 0.006737671714327339, 0.006737671714327339, 3, 1, 5.511188000000001, 16.533564, 0.370551491486672, 0.3705514914866722, 5.764849613013334, 51.88364651712001, 0.002496654302891558, 0.002496654302891559),
  @(12, "sCs", "Il1b", "Il1rap", "M2", 1, 0.3333333333333333, 1.046026666666667, 3.13808, 0.006737671714327339, 0.006737671714327339, 3, 1, 2.28691, 6.860729999999999, 0.1537632015811809, 0.1537632015811809, 2.392168844266666, 21.5295195984, 0.001036005973997935, 0.001036005973997935),
  @(13, "sCs", "Il1b", "Il1rap", "sCs", 1, 0.3333333333333333, 1.046026666666667, 3.13808, 0.006737671714327339, 0.006737671714327339, 3, 1, 4.422738333333334, 13.268215, 0.2973682418150034, 0.2973682418150034, 4.626302236355556, 41.63672012720001, 0.002003569591616201, 0.002003569591616201)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($row in $rows) {
  $r = $row[0]
  for ($i = 1; $i -lt $row.Length; $i++) {
    $col = $cols[$i - 1]
    $ws.Range("$col$r").Value = $row[$i]
  }
}